$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.606.02"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "1.788.88"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'223.64"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").Value = "'0.559"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'32.85"
$ws.Range("E8").Value = "  +6.84%  "
$ws.Range("D9").Value = "'0.281"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").Value = "'0.0680"
$ws.Range("E10").Value = "  +2.82%  "
$ws.Range("D11").Value = "'0.0936"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "2.046.70"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "'11.12"
$ws.Range("E13").Value = "  +11.07%  "
$ws.Range("D14").Value = "1.796.01"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "34.603.62"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.632"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "'4.29"
$ws.Range("E17").Value = "  +2.55%  "
$ws.Range("D18").Value = "'68.56"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "'253.58"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").Value = "0.0₃0774"
$ws.Range("E20").Value = "  +4.86%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "'10.44"
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("D23").Value = "'4.23"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'158.38"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").Value = "'16.33"
$ws.Range("E26").Value = "  -1.02%  "
$ws.Range("D27").Value = "'7.10"
$ws.Range("E27").Value = "  +2.13%  "
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.0515"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'3.75"
$ws.Range("E31").Value = "  -1.58%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "'3.57"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "'1.85"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("D35").Value = "1.443.43"
$ws.Range("E35").Value = "  -2.81%  "
$ws.Range("D36").Value = "'1.06"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.632"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0189"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("D39").Value = "'83.14"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("E40").Value = "  +4.84%  "
$ws.Range("D41").Value = "'2.36"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").Value = "'0.900"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").Value = "'0.0503"
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("E46").Value = "  -2.23%  "
$ws.Range("D47").Value = "1.944.57"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").Value = "'104.57"
$ws.Range("E48").Value = "  +7.00%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'11.97"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").Value = "'49.27"
$ws.Range("E51").Value = "  -2.96%  "
